$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Append " San Pablo" after "Universidad Catolica Boliviana" as two new
#    runs (one holding a single space, one holding "San Pablo"), each with
#    the same Helvetica Neue / 12pt / non-bold / non-italic formatting used
#    elsewhere in the byline.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$found = $anchor.Find.Execute("Universidad Catolica Boliviana", $true, $false, `
    $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $anchor.Collapse(0)
    $insertStart = $anchor.Start

    # Insert " San Pablo" in one shot so the new text first inherits the
    # formatting of the run it is attached to (keeps every rPr toggle that
    # was already explicit on the neighboring run), then bump the size to
    # 12pt (sz/szCs 24) to match the rest of the byline's larger font.
    $anchor.InsertAfter(" San Pablo")
    $wholeInsert = $d.Range($insertStart, $insertStart + 10)
    $wholeInsert.Font.Size = 12

    # Split "San Pablo" off into its own run (distinct from the leading
    # space run) by touching one of its formatting toggles; re-applying the
    # same effective value keeps the visible formatting unchanged.
    $sanPabloRange = $d.Range($insertStart + 1, $insertStart + 10)
    $sanPabloRange.Font.Italic = 1
    $sanPabloRange.Font.Italic = 0
}

# ---------------------------------------------------------------------------
# 2) Tighten the line spacing (auto, 206 -> 204 twentieths) on the "Resumen"
#    paragraph and the "Indice de terminos" paragraph.
# ---------------------------------------------------------------------------
$resumenRange = $d.Content
$foundResumen = $resumenRange.Find.Execute("Resumen", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)
if ($foundResumen) {
    $resumenPara = $resumenRange.Paragraphs.Item(1)
    $resumenPara.Format.LineSpacingRule = 5
    $resumenPara.Format.LineSpacing = 10.2
}

$indiceRange = $d.Content
$foundIndice = $indiceRange.Find.Execute("ndice de t", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0)
if ($foundIndice) {
    $indicePara = $indiceRange.Paragraphs.Item(1)
    $indicePara.Format.LineSpacingRule = 5
    $indicePara.Format.LineSpacing = 10.2
}
